$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C7) from 45207 (2023-10-08)
# to 45208 (2023-10-09) for each of the six data rows.
foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45208
}
